# ============================================================================
# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and
# excel sheets
#
# 1) Clean up "ODI Batting Extra": remove placeholder empty-string cells so
#    they become genuinely blank (matches upstream scraper now skipping
#    empty attributes instead of writing them out as "").
# 2) Add a new "ODI Bowling Extra" sheet (mirrors "ODI Batting Extra") with
#    MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Part 1: "ODI Batting Extra" -- blank out the empty placeholder cells
# ----------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$cellsToBlank = @(
    "C2", "D2", "E2",
    "C4", "D4", "E4",
    "C5", "D5", "E5",
    "B6", "C6", "D6", "E6",
    "E7",
    "B8", "C8", "D8", "E8",
    "C10", "D10", "E10",
    "B11", "C11", "D11", "E11",
    "B12", "C12", "D12", "E12",
    "B13", "C13", "D13", "E13",
    "B14", "C14", "D14", "E14",
    "C15", "D15", "E15",
    "B16", "C16", "D16", "E16",
    "B17", "C17", "D17", "E17",
    "B18", "C18", "D18", "E18", "F18",
    "B19", "C19", "D19", "E19", "F19",
    "B20", "C20", "D20", "E20", "F20",
    "B21", "C21", "D21", "E21", "F21"
)

foreach ($addr in $cellsToBlank) {
    $battingExtra.Range($addr).ClearContents()
}

# ----------------------------------------------------------------------
# Part 2: add the new "ODI Bowling Extra" sheet at the end of the workbook
# ----------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row (bold / centered / bordered, matching the other "Extra" sheet)
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $bowlingExtra.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: row, MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
$data = @(
    @(2, "3519", $null, $null),
    @(3, "3521", "0", "20.00%"),
    @(4, "3524", "0", $null),
    @(5, "3525", "1", "20.00%"),
    @(6, "3529", "0", "20.00%"),
    @(7, "3530", "1", "20.00%"),
    @(8, "3562", $null, $null),
    @(9, "3563", "1", $null),
    @(10, "3564", "1", "10.00%"),
    @(11, "3587", "0", $null),
    @(12, "3588", $null, $null),
    @(13, "3601", "0", "10.00%"),
    @(14, "3603", "0", "10.00%"),
    @(15, "3688", "1", "40.00%"),
    @(16, "3689", $null, $null),
    @(17, "3692", $null, $null),
    @(18, "3875", $null, $null),
    @(19, "3876", "0", "20.00%"),
    @(20, "3877", $null, $null),
    @(21, "3878", $null, $null)
)

foreach ($row in $data) {
    $rowIndex = $row[0]

    $cellA = $bowlingExtra.Cells.Item($rowIndex, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[1]

    $cellB = $bowlingExtra.Cells.Item($rowIndex, 2)
    $cellB.NumberFormat = "@"
    if ($row[2] -ne $null) {
        $cellB.Value = $row[2]
    }

    $cellC = $bowlingExtra.Cells.Item($rowIndex, 3)
    $cellC.NumberFormat = "@"
    if ($row[3] -ne $null) {
        $cellC.Value = $row[3]
    }
}

$bowlingExtra.Range("A1").Select() | Out-Null
